# Applies the "Resumo de Inscrições - Superior" update:
# A handful of courses/campi had their Inscritos (E), Pagos (F) and/or
# Inscrições homologadas (H) counts bumped up (new inscriptions/payments
# registered). Isenções deferidas (G) is unaffected (stays 0 everywhere).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value = 11
$ws.Range("H3").Value = 11

$ws.Range("F4").Value = 5
$ws.Range("H4").Value = 5

$ws.Range("F9").Value = 7
$ws.Range("H9").Value = 7

$ws.Range("F15").Value = 73
$ws.Range("H15").Value = 73

$ws.Range("E17").Value = 93
$ws.Range("F17").Value = 41
$ws.Range("H17").Value = 41

$ws.Range("E18").Value = 96
$ws.Range("F18").Value = 37
$ws.Range("H18").Value = 37

$ws.Range("E19").Value = 39

$ws.Range("E28").Value = 12

$ws.Range("E35").Value = 5

$ws.Range("F36").Value = 33
$ws.Range("H36").Value = 33

$ws.Range("F37").Value = 20
$ws.Range("H37").Value = 20

$ws.Range("E38").Value = 56

$ws.Range("F44").Value = 10
$ws.Range("H44").Value = 10

$ws.Range("E47").Value = 49
$ws.Range("F47").Value = 30
$ws.Range("H47").Value = 30

$ws.Range("F51").Value = 5
$ws.Range("H51").Value = 5

$ws.Range("E59").Value = 7
$ws.Range("F59").Value = 3
$ws.Range("H59").Value = 3

$ws.Range("E61").Value = 25

$ws.Range("F62").Value = 10
$ws.Range("H62").Value = 10

$ws.Range("E63").Value = 26

$ws.Range("E70").Value = 36
$ws.Range("F70").Value = 16
$ws.Range("H70").Value = 16

$ws.Range("E73").Value = 25

$ws.Range("F76").Value = 16
$ws.Range("H76").Value = 16

$ws.Range("E78").Value = 40
$ws.Range("F78").Value = 15
$ws.Range("H78").Value = 15

$ws.Range("E81").Value = 11

$wb.Save()
